# "fix error export data"
# The product export sheet previously contained 4 stale data rows (rows 2-5).
# This fixes the export so it only contains the current, correct product
# record in row 2, with refreshed values, and removes the leftover rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale/incorrect extra data rows (old rows 3, 4 and 5).
$ws.Rows("3:5").Delete()

# Refresh the remaining data row with the corrected export values.
$ws.Range("B2").Value = "2023-12-29 04:03:55"
$ws.Range("C2").Value = "Pupuk"
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = "0,00"
